$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "25.819.52"
Set-TextValue "E2" "  -0.54%  "
Set-TextValue "D3" "1.632.94"
Set-TextValue "E3" "  -0.46%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "215.68"
Set-TextValue "E5" "  +0.39%  "
Set-TextValue "D6" "0.5067"
Set-TextValue "E6" "  -0.17%  "
Set-TextValue "D7" "1.002"
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "E8" "  -0.14%  "
Set-TextValue "D9" "0.06424"
Set-TextValue "E9" "  +0.98%  "
Set-TextValue "D10" "19.49"
Set-TextValue "E10" "  -1.93%  "
Set-TextValue "D11" "0.07801"
Set-TextValue "E11" "  +0.73%  "
Set-TextValue "D12" "4.279"
Set-TextValue "E12" "  -0.43%  "
Set-TextValue "B13" "WrappedliquidstakedEther2.0"
Set-TextValue "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "1.856.68"
Set-TextValue "E13" "  -0.59%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.631.45"
Set-TextValue "E14" "  -0.62%  "
Set-TextValue "D15" "0.5599"
Set-TextValue "E15" "  +2.37%  "
Set-TextValue "B16" "Litecoin"
Set-TextValue "C16" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "63.08"
Set-TextValue "E16" "  -1.82%  "
Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0₅7582"
Set-TextValue "E17" "  -2.19%  "
Set-TextValue "D18" "25.822.64"
Set-TextValue "E18" "  -0.67%  "
Set-TextValue "D19" "1.003"
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "D20" "194.62"
Set-TextValue "E20" "  -0.94%  "
Set-TextValue "D21" "4.326"
Set-TextValue "E21" "  -3.08%  "
Set-TextValue "D22" "9.846"
Set-TextValue "E22" "  -1.19%  "
Set-TextValue "D23" "6.046"
Set-TextValue "E23" "  -1.64%  "
Set-TextValue "E24" "  +0.00%  "
Set-TextValue "D25" "1.796"
Set-TextValue "E25" "  -4.97%  "
Set-TextValue "E26" "  +0.53%  "
Set-TextValue "D27" "140.34"
Set-TextValue "E27" "  -1.59%  "
Set-TextValue "D28" "6.751"
Set-TextValue "E28" "  -1.69%  "
Set-TextValue "D29" "15.45"
Set-TextValue "E29" "  -1.23%  "
Set-TextValue "D30" "1.238"
Set-TextValue "E30" "  -0.11%  "
Set-TextValue "D31" "0.04882"
Set-TextValue "E31" "  -0.31%  "
Set-TextValue "E32" "  +0.50%  "
Set-TextValue "D33" "3.216"
Set-TextValue "E33" "  +0.37%  "
Set-TextValue "D34" "1.551"
Set-TextValue "E34" "  +0.16%  "
Set-TextValue "D35" "2.378"
Set-TextValue "E35" "  +0.07%  "
Set-TextValue "D36" "0.8978"
Set-TextValue "E36" "  -2.26%  "
Set-TextValue "D37" "2.568"
Set-TextValue "E37" "  +0.05%  "
Set-TextValue "D38" "1.128.57"
Set-TextValue "E39" "  -0.85%  "
Set-TextValue "E40" "  -0.78%  "
Set-TextValue "D41" "0.9941"
Set-TextValue "E41" "  -0.75%  "
Set-TextValue "D42" "5.530"
Set-TextValue "E42" "  -1.19%  "
Set-TextValue "D43" "0.7984"
Set-TextValue "E43" "  -0.62%  "
Set-TextValue "D44" "97.29"
Set-TextValue "E44" "  -1.46%  "
Set-TextValue "D45" "1.783.29"
Set-TextValue "E45" "  +0.29%  "
Set-TextValue "D46" "0.0₈111"
Set-TextValue "E46" "  -6.74%  "
Set-TextValue "D47" "0.4439"
Set-TextValue "E47" "  -1.90%  "
Set-TextValue "D48" "55.35"
Set-TextValue "E48" "  +0.11%  "
Set-TextValue "E49" "  -2.67%  "
Set-TextValue "D50" "7.650"
Set-TextValue "E50" "  +2.02%  "
Set-TextValue "D51" "0.9987"
Set-TextValue "E51" "  -0.10%  "
